$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet is an alphabetically-sorted English -> translation ("TBT"
# placeholder) list. Four new English phrases were added to the list
# (each already correctly alphabetically positioned), pushing every row
# below each insertion point down by one. Applying the inserts in
# ascending row order (each Insert+Value pair executed immediately,
# before moving on to the next) lets each subsequent target row number
# be used as-is, since earlier inserts have already shifted the sheet by
# the time we reach it.

# New row 33: "Clinical data not provided"
$ws.Rows(33).Insert()
$ws.Range("A33").Value = "Clinical data not provided"
$ws.Range("B33").Value = "TBT"

# New row 80: "Lab data not provided"
$ws.Rows(80).Insert()
$ws.Range("A80").Value = "Lab data not provided"
$ws.Range("B80").Value = "TBT"

# New row 95: "No .acorn has been generated"
$ws.Rows(95).Insert()
$ws.Range("A95").Value = "No .acorn has been generated"
$ws.Range("B95").Value = "TBT"

# New row 96: "No .acorn has been saved"
$ws.Rows(96).Insert()
$ws.Range("A96").Value = "No .acorn has been saved"
$ws.Range("B96").Value = "TBT"
